# Update the build-version timestamp throughout the workbook.
$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on January 30 2026 16.19.47 EST)"
$newVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"

$newCitation = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Flying Eagle Mine, United States, M3466, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "About" sheet ---
$aboutSheet = $wb.Worksheets.Item("About")
$aboutSheet.Range("A2").Value = "Version: $newVersion"
$aboutSheet.Range("A6").Value = $newCitation

# --- "Boundaries and methane sources" sheet ---
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")
for ($row = 2; $row -le 7; $row++) {
    $cell = $dataSheet.Cells.Item($row, 19)  # Column S = build_version
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
